$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each line: Row,D(Fecha serial),J(Volumen),K(Precio minimo),L(Precio maximo),M(Precio promedio ponderado),N(Unidad comercializacion),O(Origen),P(Precio $/Kg),Q(Kg o Unidades)
$data = @(
"271,45027,20,6000,6000,6000,`$/paquete,Región Metropolitana,6000,1",
"272,44971,50,5000,5000,5000,`$/paquete,Región de La Araucanía,5000,1",
"273,44971,80,4000,4000,4000,`$/paquete,Región del Maule,4000,1",
"274,44553,65,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"275,44985,40,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"276,44222,155,6000,6000,6000,`$/paquete,Región de La Araucanía,6000,1",
"277,44518,95,8000,8000,8000,`$/paquete,Región del Maule,8000,1",
"278,44420,95,8000,9000,8474,`$/paquete,Región de Arica y Parinacota,8474,1",
"279,44453,20,8000,8000,8000,`$/paquete,Región de Arica y Parinacota,8000,1",
"280,44757,40,6000,6000,6000,`$/paquete,Región de Arica y Parinacota,6000,1",
"281,44648,30,5500,5500,5500,`$/paquete,Región del Maule,5500,1",
"282,44809,80,6000,6000,6000,`$/paquete,Región de Arica y Parinacota,6000,1",
"283,44634,50,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"284,44953,55,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"285,44616,40,7000,7000,7000,`$/paquete,Región de La Araucanía,7000,1",
"286,44187,40,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"287,44418,45,6000,6000,6000,`$/paquete,Región de Arica y Parinacota,6000,1",
"288,44413,30,8000,8000,8000,`$/paquete,Región de Arica y Parinacota,8000,1",
"289,44266,65,5000,6000,5462,`$/paquete,Región de La Araucanía,5462,1",
"290,44382,80,5000,5000,5000,`$/paquete,Región de Arica y Parinacota,5000,1",
"291,44907,80,8000,8000,8000,`$/paquete,Región del Maule,8000,1",
"292,44278,45,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"293,44286,70,6000,7000,6429,`$/paquete,Región de La Araucanía,6429,1",
"294,44950,25,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"295,44362,25,5000,5000,5000,`$/paquete,Región de Arica y Parinacota,5000,1",
"296,44988,50,5000,6000,5600,`$/paquete,Región de La Araucanía,5600,1",
"297,44988,40,5000,6000,5500,`$/paquete,Región del Maule,5500,1",
"298,45006,25,5000,5000,5000,`$/paquete,Región de La Araucanía,5000,1",
"299,44181,35,5000,5000,5000,`$/paquete,Región de Arica y Parinacota,5000,1",
"300,44679,50,5000,5000,5000,`$/paquete,Región de Arica y Parinacota,5000,1",
"301,44627,15,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"302,44544,75,6000,7000,6467,`$/paquete,Región del Maule,6467,1",
"303,44677,30,6000,6000,6000,`$/paquete,Región de Arica y Parinacota,6000,1",
"304,44259,40,5000,5000,5000,`$/paquete,Región de La Araucanía,5000,1",
"305,44259,70,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"306,44340,50,5000,5000,5000,`$/paquete,Región de Arica y Parinacota,5000,1",
"307,44188,40,6000,6000,6000,`$/docena,Región de La Araucanía,5000,1.2",
"308,44188,80,4000,4000,4000,`$/docena,Región del Maule,3333,1.2",
"309,44754,30,6000,6000,6000,`$/paquete,Región de Arica y Parinacota,6000,1",
"310,44200,40,6000,6000,6000,`$/docena,Región de La Araucanía,5000,1.2",
"311,44200,50,5000,5000,5000,`$/docena,Región del Maule,4167,1.2",
"312,44603,35,6000,6000,6000,`$/paquete,Región de La Araucanía,6000,1",
"313,44767,100,5000,5000,5000,`$/paquete,Región de Arica y Parinacota,5000,1",
"314,44760,65,5500,7000,6077,`$/paquete,Región de Arica y Parinacota,6077,1",
"315,44571,65,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"316,44595,40,7000,7000,7000,`$/paquete,Región de La Araucanía,7000,1",
"317,44595,100,7000,7000,7000,`$/paquete,Región del Maule,7000,1",
"318,44662,80,6000,6000,6000,`$/paquete,Región Metropolitana,6000,1",
"319,44189,40,5000,6000,5500,`$/docena,Región de La Araucanía,4583,1.2",
"320,44189,40,3000,3000,3000,`$/paquete,Región de Arica y Parinacota,3000,1",
"321,44385,50,5000,5000,5000,`$/paquete,Región de Arica y Parinacota,5000,1",
"322,44371,60,5000,6000,5500,`$/paquete,Región de Arica y Parinacota,5500,1",
"323,44364,45,5000,5000,5000,`$/paquete,Región de Arica y Parinacota,5000,1",
"324,45008,35,5000,5000,5000,`$/paquete,Región de La Araucanía,5000,1",
"325,45008,65,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"326,44873,25,10000,10000,10000,`$/paquete,Región de Arica y Parinacota,10000,1",
"327,44936,65,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"328,44473,40,7000,7000,7000,`$/paquete,Región de Arica y Parinacota,7000,1",
"329,44951,25,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"330,44218,30,5000,5000,5000,`$/paquete,Región de La Araucanía,5000,1",
"331,44218,40,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"332,44910,140,8000,8000,8000,`$/paquete,Región Metropolitana,8000,1",
"333,44910,80,8000,8000,8000,`$/paquete,Región del Maule,8000,1",
"334,44879,100,9000,9000,9000,`$/paquete,Región Metropolitana,9000,1",
"335,44239,90,5000,6000,5611,`$/paquete,Región de La Araucanía,5611,1",
"336,44239,85,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"337,44883,50,8000,8000,8000,`$/paquete,Región Metropolitana,8000,1",
"338,44232,45,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"339,44868,180,9000,10000,9500,`$/paquete,Región de Arica y Parinacota,9500,1",
"340,44847,65,9000,9000,9000,`$/paquete,Región de Arica y Parinacota,9000,1",
"341,44921,65,8000,8000,8000,`$/paquete,Región Metropolitana,8000,1",
"342,44987,80,6000,6000,6000,`$/paquete,Región de La Araucanía,6000,1",
"343,44987,100,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"344,44257,30,5000,5000,5000,`$/paquete,Región de La Araucanía,5000,1",
"345,44257,80,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"346,44264,40,5000,5000,5000,`$/paquete,Región de La Araucanía,5000,1",
"347,44264,35,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"348,44901,35,8000,8000,8000,`$/paquete,Región Metropolitana,8000,1",
"349,44176,10,10000,10000,10000,`$/docena,Región de La Araucanía,8333,1.2",
"350,44176,20,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"351,44474,20,7000,7000,7000,`$/atado,Región de Arica y Parinacota,7000,1",
"352,45002,100,4000,5000,4700,`$/paquete,Región de La Araucanía,4700,1",
"353,45002,40,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"354,45012,50,5000,5000,5000,`$/paquete,Región de La Araucanía,5000,1",
"355,44705,30,5000,5000,5000,`$/paquete,Región de Arica y Parinacota,5000,1",
"356,44448,40,7000,7000,7000,`$/paquete,Región de Arica y Parinacota,7000,1",
"357,44970,50,6000,6000,6000,`$/paquete,Región de La Araucanía,6000,1",
"358,44970,100,4000,4000,4000,`$/paquete,Región del Maule,4000,1",
"359,45015,80,5000,6000,5500,`$/paquete,Región de La Araucanía,5500,1",
"360,44186,50,8000,8000,8000,`$/docena,Región de Arica y Parinacota,6667,1.2",
"361,44186,50,3000,3000,3000,`$/paquete,Región de Arica y Parinacota,3000,1",
"362,44179,100,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"363,44172,60,4000,5000,4583,`$/paquete,Región del Maule,4583,1",
"364,44550,125,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"365,44952,45,4000,4000,4000,`$/paquete,Región de La Araucanía,4000,1",
"366,44952,55,6000,6000,6000,`$/paquete,Región del Maule,6000,1",
"367,44252,95,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"368,44855,20,8000,8000,8000,`$/paquete,Región de Arica y Parinacota,8000,1",
"369,44168,65,4000,4000,4000,`$/paquete,Región del Maule,4000,1",
"370,44875,65,10000,10000,10000,`$/paquete,Región de Arica y Parinacota,10000,1",
"371,44642,25,5000,5000,5000,`$/paquete,Región del Maule,5000,1",
"372,44454,20,8000,8000,8000,`$/paquete,Región de Arica y Parinacota,8000,1",
"373,44426,40,6000,6000,6000,`$/paquete,Región de Arica y Parinacota,6000,1",
"374,44606,40,7000,7000,7000,`$/paquete,Región del Maule,7000,1",
"375,44526,20,3500,3500,3500,`$/paquete,Región del Maule,3500,1",
"376,44918,125,8000,8000,8000,`$/paquete,Región Metropolitana,8000,1",
"377,44243,75,5000,6000,5400,`$/paquete,Región de La Araucanía,5400,1",
"378,44217,50,5000,5000,5000,`$/paquete,Región de La Araucanía,5000,1",
"379,44217,80,5000,5000,5000,`$/paquete,Región del Maule,5000,1"
)

foreach ($line in $data) {
    $f = $line -split ","
    $r = [int]$f[0]

    $ws.Cells.Item($r, 4).Value = [double]$f[1]
    $ws.Cells.Item($r, 10).Value = [double]$f[2]
    $ws.Cells.Item($r, 11).Value = [double]$f[3]
    $ws.Cells.Item($r, 12).Value = [double]$f[4]
    $ws.Cells.Item($r, 13).Value = [double]$f[5]
    $ws.Cells.Item($r, 14).Value = $f[6]
    $ws.Cells.Item($r, 15).Value = $f[7]
    $ws.Cells.Item($r, 16).Value = [double]$f[8]
    $ws.Cells.Item($r, 17).Value = [double]$f[9]
}

# Row 379 is brand new: fill in the constant columns that are the same for every
# row in this product block (A,B,C,E,F,G,H,I,R).
$ws.Cells.Item(379, 1).Value = 10
$ws.Cells.Item(379, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(379, 3).Value = "La Araucanía"
$ws.Cells.Item(379, 5).Value = 9
$ws.Cells.Item(379, 6).Value = 100112052
$ws.Cells.Item(379, 7).Value = "Albahaca"
$ws.Cells.Item(379, 8).Value = "Sin especificar"
$ws.Cells.Item(379, 9).Value = "Primera"
$ws.Cells.Item(379, 18).Value = "Hortaliza"

# D379 needs the date style (s="2") like the rest of column D; copy the format
# from D378 (value for D379 was already set in the loop above).
$ws.Range("D378").Copy()
$ws.Range("D379").PasteSpecial(-4122)  # xlPasteFormats
